$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 177 - add MOVE value and update VIX value
$ws.Range("B177").Value = 53.1
$ws.Range("C177").Value = 19.37

# Add new rows 178-182 with new dates and data
$ws.Range("A178").Value = "14-09-2021"
$ws.Range("B178").Value = 51.73
$ws.Range("C178").Value = 19.46

$ws.Range("A179").Value = "15-09-2021"
$ws.Range("B179").Value = 53.93
$ws.Range("C179").Value = 18.18

$ws.Range("A180").Value = "16-09-2021"
$ws.Range("B180").Value = 53.7
$ws.Range("C180").Value = 18.69

$ws.Range("A181").Value = "17-09-2021"
$ws.Range("B181").Value = 56.06
$ws.Range("C181").Value = 20.81

$ws.Range("A182").Value = "20-09-2021"
$ws.Range("C182").Value = 26.23
